# [Refactor] 변수명 통일 중 - name을 displayName/EnemyName으로 통일, player를 character로 통일
#
# For this sheet (characterCardLevelInfo_data), the header cell that used
# to hold the literal field name "name" is renamed to "display_name".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "name" column header to "display_name"
$ws.Range("A1").Value = "display_name"

# The header text got longer, so column A is re-sized to fit its content
# again (it was a bestFit/autofit column before the rename).
$ws.Columns.Item(1).AutoFit() | Out-Null

# Leave the cursor/selection on a single cell (E4), matching the saved
# workbook state after the edit.
$ws.Range("E4").Select() | Out-Null
